$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("P2").Value = 84
$ws.Range("Q2").Value = 42
$ws.Range("R2").Formula = "=P2+Q2"

# Row 3
$ws.Range("P3").Formula = "=P2+Q3"
$ws.Range("Q3").Value = 7
$ws.Range("R3").Formula = "=Q3"

# Row 4
$ws.Range("P4").Formula = "=P3+Q4"
$ws.Range("Q4").Value = 7
$ws.Range("R4").Formula = "=R3+Q4"

# Row 5
$ws.Range("P5").Formula = "=P4+Q5"
$ws.Range("Q5").Value = 7
$ws.Range("R5").Formula = "=R4+Q5"

# Row 6
$ws.Range("P6").Formula = "=P5+Q6"
$ws.Range("Q6").Value = 7
$ws.Range("R6").Formula = "=R5+Q6"

# Row 7
$ws.Range("P7").Formula = "=P6+Q7"
$ws.Range("Q7").Value = 7
$ws.Range("R7").Formula = "=R6+Q7"

# Row 8 (no P8)
$ws.Range("Q8").Value = 7
$ws.Range("R8").Formula = "=R7+Q8"

# Row 9 (no P9, no R9)
$ws.Range("Q9").Formula = "=SUM(Q3:Q8)"

# Update selection to match the target view state
$ws.Range("R6:R8").Select()
